$d = $word.ActiveDocument

# wdAlignParagraphJustify = 3, wdAlignParagraphLeft = 0
# The three trailing paragraphs carry an explicit "justify" (w:jc w:val="both")
# paragraph property; the edit clears it so the paragraphs fall back to the
# document default (no explicit alignment), matching the other paragraphs
# in the document which already have an empty <w:pPr/>.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Format.Alignment -eq 3) {
        $p.Format.Alignment = 0
    }
}
